$d = $word.ActiveDocument

# --- Locate the sentence that gets extended -------------------------------
# "Check which part takes most time and t" is the run immediately before the
# "_GoBack" bookmark; collapse the found range to its end so we insert right
# at the bookmark's (current) position.
$anchor = $d.Content
$found = $anchor.Find.Execute("Check which part takes most time and t", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor text 'Check which part takes most time and t'"
}
$anchor.Collapse(0)  # wdCollapseEnd

# --- Insert "ry to improve. (possible changes at " before the bookmark ----
# The bookmark is a collapsed range sitting right here; inserting text via
# InsertAfter on a range that ends exactly at the bookmark leaves the
# bookmark anchored after the newly inserted text (mirrors the diff, where
# the bookmark stays put while new runs are added in front of it).
$anchor.InsertAfter("ry to improve. (possible changes at ")

# --- Replace the old "ry to improve." run (now after the bookmark) --------
# with "line 390)".
$bm = $d.Bookmarks.Item("_GoBack")
$tail = $d.Range($bm.End, $bm.End)
$tailFound = $tail.Find.Execute("ry to improve.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $tailFound) {
    throw "Could not find trailing text 'ry to improve.' after the bookmark"
}
$tail.Text = "line 390)"
